$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text like "1.00" / "584.17" that Excel would
# otherwise auto-convert to numbers; force them through a Text number format
# and restore the default "Normal" style afterwards so the stored cell keeps
# no explicit style (matching the original inlineStr text cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "63.543.11"
$ws.Range("E2").Value = "  +0.68%  "
Set-TextValue $ws.Range("D3") "3.099.49"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "584.17"
$ws.Range("E5").Value = "  -0.15%  "
Set-TextValue $ws.Range("D6") "144.97"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.08%  "
Set-TextValue $ws.Range("D8") "3.092.57"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  +7.28%  "
$ws.Range("E11").Value = "  -2.89%  "
Set-TextValue $ws.Range("D12") "0.456"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("E13").Value = "  -0.74%  "
Set-TextValue $ws.Range("D14") "37.06"
$ws.Range("E14").Value = "  +4.31%  "
$ws.Range("E15").Value = "  -1.17%  "
Set-TextValue $ws.Range("D16") "3.612.29"
$ws.Range("E16").Value = "  -0.41%  "
Set-TextValue $ws.Range("D17") "63.376.52"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  -1.45%  "
Set-TextValue $ws.Range("D19") "3.088.02"
$ws.Range("E19").Value = "  -0.83%  "
Set-TextValue $ws.Range("D20") "460.29"
$ws.Range("E20").Value = "  -1.34%  "
Set-TextValue $ws.Range("D21") "14.21"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("E24").Value = "  -2.81%  "
Set-TextValue $ws.Range("D25") "81.05"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  +2.91%  "
$ws.Range("E27").Value = "  +0.03%  "
Set-TextValue $ws.Range("D28") "9.19"
$ws.Range("E28").Value = "  +10.00%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D29") "2.67"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  -0.07%  "
Set-TextValue $ws.Range("D31") "2.20"
$ws.Range("E31").Value = "  -1.49%  "
Set-TextValue $ws.Range("D32") "6.97"
$ws.Range("E32").Value = "  +1.99%  "
Set-TextValue $ws.Range("D33") "0.111"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  -0.95%  "
Set-TextValue $ws.Range("D35") "0.0₃0847"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("E36").Value = "  -0.61%  "
Set-TextValue $ws.Range("D37") "3.37"
$ws.Range("E37").Value = "  +2.49%  "
Set-TextValue $ws.Range("D38") "2.30"
$ws.Range("E38").Value = "  -4.83%  "
$ws.Range("E39").Value = "  -0.50%  "
Set-TextValue $ws.Range("D40") "50.23"
$ws.Range("E40").Value = "  -1.40%  "
Set-TextValue $ws.Range("D41") "435.82"
$ws.Range("E41").Value = "  +1.25%  "
Set-TextValue $ws.Range("D42") "8.68"
$ws.Range("E42").Value = "  -0.43%  "
Set-TextValue $ws.Range("D43") "0.0368"
$ws.Range("E43").Value = "  -0.31%  "
Set-TextValue $ws.Range("D44") "2.883.05"
$ws.Range("E44").Value = "  -1.50%  "
Set-TextValue $ws.Range("D45") "0.274"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("E46").Value = "  -3.00%  "
Set-TextValue $ws.Range("D47") "36.45"
$ws.Range("E47").Value = "  +3.04%  "
Set-TextValue $ws.Range("D48") "124.99"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -1.18%  "
Set-TextValue $ws.Range("D51") "24.07"
$ws.Range("E51").Value = "  -1.77%  "
